$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values to match repulled data
$ws.Range("F2").Value = -9
$ws.Range("F4").Value = 8
$ws.Range("F6").Value = -1
$ws.Range("F8").Value = 2
$ws.Range("F9").Value = 1
$ws.Range("F10").Value = 0
$ws.Range("F11").Value = 1
